{"js": "// The edit swaps the word order in the Finnish header/footer line from\n// \"havainnointijaksot vuonna Saappaiden t\u00e4hdist\u00f6 2022: ...\" to\n// \"Saappaiden t\u00e4hdist\u00f6 havainnointijaksot vuonna 2022: ...\" everywhere it\n// occurs in the document body (4 occurrences).\n\nconst oldText = \"havainnointijaksot vuonna Saappaiden t\u00e4hdist\u00f6 2022\";\nconst newText = \"Saappaiden t\u00e4hdist\u00f6 havainnointijaksot vuonna 2022\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The edit swaps the word order in the Finnish header/footer line from\n# \"havainnointijaksot vuonna Saappaiden t\u00e4hdist\u00f6 2022: ...\" to\n# \"Saappaiden t\u00e4hdist\u00f6 havainnointijaksot vuonna 2022: ...\" everywhere it\n# occurs in the document body (4 occurrences).\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"havainnointijaksot vuonna Saappaiden t\u00e4hdist\u00f6 2022\"\n$find.Replacement.Text = \"Saappaiden t\u00e4hdist\u00f6 havainnointijaksot vuonna 2022\"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, [ref]$find.MatchWholeWord, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$find.Replacement.Text, [ref]2) | Out-Null\n"}
